$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.222.33'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '1.565.41'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.84'
$ws.Range("E5").Value = '  +1.51%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.14'
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("E10").Value = '  -0.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0873'
$ws.Range("E11").Value = '  +2.36%  '
$ws.Range("D12").Value = '1.788.41'
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").Value = '1.552.98'
$ws.Range("E13").Value = '  -0.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.75'
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("E15").Value = '  -0.55%  '
$ws.Range("D16").Value = '27.193.18'
$ws.Range("E16").Value = '  +0.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.87'
$ws.Range("E17").Value = '  -0.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '217.75'
$ws.Range("E18").Value = '  +0.72%  '
$ws.Range("E19").Value = '  +1.17%  '
$ws.Range("D20").Value = '0.0₃0702'
$ws.Range("E20").Value = '  -0.73%  '
$ws.Range("E22").Value = '  +0.49%  '
$ws.Range("E23").Value = '  +1.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.65'
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.62'
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  +1.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.02'
$ws.Range("E28").Value = '  -0.47%  '
$ws.Range("E30").Value = '  +1.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0470'
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("D33").Value = '1.460.49'
$ws.Range("E33").Value = '  +1.93%  '
$ws.Range("E34").Value = '  +0.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.11'
$ws.Range("E35").Value = '  +4.34%  '
$ws.Range("E36").Value = '  +1.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.35'
$ws.Range("E37").Value = '  +0.81%  '
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("E39").Value = '  +1.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.87'
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.813'
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("E43").Value = '  +1.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.983'
$ws.Range("E44").Value = '  -1.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.42'
$ws.Range("E45").Value = '  -0.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.76'
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("D47").Value = '1.700.25'
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.85'
$ws.Range("E48").Value = '  -1.66%  '
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0525'
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0947'
$ws.Range("E51").Value = '  -1.40%  '
